$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 287 (shifts existing rows 287..396 down to 288..397)
$ws.Rows(287).Insert()

# Fill the newly inserted row 287 with the new record's data
$ws.Range("A287").Value = 5
$ws.Range("B287").Value = "Macroferia Regional de Talca"
$ws.Range("C287").Value = "Maule"
$ws.Range("D287").Value = 44825
$ws.Range("E287").Value = 7
$ws.Range("F287").Value = 100114014
$ws.Range("G287").Value = "Betarraga"
$ws.Range("H287").Value = "Sin especificar"
$ws.Range("I287").Value = "Segunda"
$ws.Range("J287").Value = 4000
$ws.Range("K287").Value = 800
$ws.Range("L287").Value = 800
$ws.Range("M287").Value = 800
$ws.Range("N287").Value = "$/paquete 5 unidades"
$ws.Range("O287").Value = "Región del Maule"
$ws.Range("P287").Value = 160
$ws.Range("Q287").Value = 5
$ws.Range("R287").Value = "Hortaliza"
